$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the Price (D) and/or Volume (E) columns hold numeric-looking
# text (e.g. "67.245.64", "-3.11%") that must stay text, not become a number.
$textFormatRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $textFormatRows) {
    $ws.Range("D$r`:E$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "67.245.64"
$ws.Range("E2").Value = "  -3.11%  "

# Row 3
$ws.Range("D3").Value = "3.492.50"
$ws.Range("E3").Value = "  -4.66%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "600.78"
$ws.Range("E5").Value = "  -3.65%  "

# Row 6
$ws.Range("D6").Value = "148.34"
$ws.Range("E6").Value = "  -6.83%  "

# Row 7
$ws.Range("D7").Value = "3.491.43"
$ws.Range("E7").Value = "  -4.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  -3.23%  "

# Row 10
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  -5.07%  "

# Row 11
$ws.Range("D11").Value = "6.99"
$ws.Range("E11").Value = "  -3.00%  "

# Row 12
$ws.Range("D12").Value = "0.421"
$ws.Range("E12").Value = "  -4.13%  "

# Row 13
$ws.Range("D13").Value = "0.0000215"
$ws.Range("E13").Value = "  -5.94%  "

# Row 14
$ws.Range("D14").Value = "4.083.96"
$ws.Range("E14").Value = "  -4.59%  "

# Row 15
$ws.Range("D15").Value = "31.33"
$ws.Range("E15").Value = "  -2.67%  "

# Row 16
$ws.Range("D16").Value = "3.493.63"
$ws.Range("E16").Value = "  -4.87%  "

# Row 17
$ws.Range("D17").Value = "67.226.67"
$ws.Range("E17").Value = "  -3.18%  "

# Row 18
$ws.Range("E18").Value = "  -1.05%  "

# Row 19
$ws.Range("D19").Value = "6.36"
$ws.Range("E19").Value = "  -2.00%  "

# Row 20
$ws.Range("E20").Value = "  -5.16%  "

# Row 21
$ws.Range("D21").Value = "444.82"
$ws.Range("E21").Value = "  -4.71%  "

# Row 22
$ws.Range("D22").Value = "8.98"
$ws.Range("E22").Value = "  -13.03%  "

# Row 23
$ws.Range("D23").Value = "0.617"
$ws.Range("E23").Value = "  -4.79%  "

# Row 24
$ws.Range("D24").Value = "77.10"
$ws.Range("E24").Value = "  -2.66%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.07%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.632.52"
$ws.Range("E26").Value = "  -4.65%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0000126"
$ws.Range("E27").Value = "  +3.49%  "

# Row 28
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  -10.14%  "

# Row 29
$ws.Range("D29").Value = "8.12"
$ws.Range("E29").Value = "  -5.77%  "

# Row 30
$ws.Range("D30").Value = "2.46"
$ws.Range("E30").Value = "  -5.72%  "

# Row 31
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.19%  "

# Row 32
$ws.Range("D32").Value = "1.53"
$ws.Range("E32").Value = "  -8.36%  "

# Row 33
$ws.Range("D33").Value = "0.165"
$ws.Range("E33").Value = "  +0.41%  "

# Row 34
$ws.Range("D34").Value = "25.64"
$ws.Range("E34").Value = "  -3.23%  "

# Row 35
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.479.94"
$ws.Range("E35").Value = "  -5.08%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").Value = "  -5.32%  "

# Row 37
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -7.00%  "

# Row 38
$ws.Range("D38").Value = "7.96"
$ws.Range("E38").Value = "  -3.73%  "

# Row 40
$ws.Range("E40").Value = "  -0.08%  "

# Row 41
$ws.Range("D41").Value = "176.10"
$ws.Range("E41").Value = "  -1.97%  "

# Row 42
$ws.Range("D42").Value = "2.17"
$ws.Range("E42").Value = "  -2.16%  "

# Row 43
$ws.Range("D43").Value = "0.0871"
$ws.Range("E43").Value = "  -2.68%  "

# Row 44
$ws.Range("D44").Value = "5.36"
$ws.Range("E44").Value = "  -7.57%  "

# Row 45
$ws.Range("D45").Value = "0.877"
$ws.Range("E45").Value = "  -4.86%  "

# Row 46
$ws.Range("D46").Value = "45.37"
$ws.Range("E46").Value = "  -2.84%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "28.06"
$ws.Range("E47").Value = "  -4.83%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.25"
$ws.Range("E48").Value = "  +4.10%  "

# Row 49
$ws.Range("D49").Value = "2.51"
$ws.Range("E49").Value = "  -6.90%  "

# Row 50
$ws.Range("D50").Value = "7.51"
$ws.Range("E50").Value = "  -4.02%  "

# Row 51
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -3.38%  "

# Restore default (unstyled) formatting so cells match the original plain-text layout
foreach ($r in $textFormatRows) {
    $ws.Range("D$r`:E$r").Style = "Normal"
}
